# Generate Report for Handback
#
# The handback transform for the "4c62728e-0f11-4df0-a093-6cce9146ab75.md"
# file failed for both target locales. Update the localization-status
# report: flip the row's Status from "Ready for handoff" to
# "Handback transform failed" (both on the per-locale sheets and on the
# Overview roll-up sheet, which mirrors the same status per locale column),
# fill in the Error Detail column with the diagnostic message, and widen
# the Error Detail column so the message is readable.

$wb = $excel.ActiveWorkbook

# --- Overview sheet --------------------------------------------------------
# The Overview sheet mirrors each locale's Status for this file in its
# "zh-cn" (E) and "de-de" (F) columns on row 3.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Handback transform failed"
$wsOverview.Range("F3").Value = "Handback transform failed"

# --- zh-cn sheet ---------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

# Row 3 corresponds to 4c62728e-0f11-4df0-a093-6cce9146ab75.md
$wsZhCn.Range("C3").Value = "Handback transform failed"
$wsZhCn.Range("P3").Value = "Handback file name: dlglw3jw.jig is different with handoff file name: 4c62728e-0f11-4df0-a093-6cce9146ab75.507cd666f575c3f23b884ee76c4f87bdf2951f54.zh-cn."

# Widen the "Error Detail" column (P) so the new message is visible.
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17

# --- de-de sheet -----------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

# Row 3 corresponds to 4c62728e-0f11-4df0-a093-6cce9146ab75.md
$wsDeDe.Range("C3").Value = "Handback transform failed"
$wsDeDe.Range("P3").Value = "Handback file name: dlglw3jw.jig is different with handoff file name: 4c62728e-0f11-4df0-a093-6cce9146ab75.507cd666f575c3f23b884ee76c4f87bdf2951f54.de-de."

# Widen the "Error Detail" column (P) so the new message is visible.
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
